$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 1894.8572
$ws.Range("J88").Value = 2313.8
$ws.Range("L88").Value = 2313.8
$ws.Range("N88").Value = -3125.8
$ws.Range("H91").Value = 1894.8572
$ws.Range("J91").Value = 2313.8
$ws.Range("L91").Value = 2313.8
$ws.Range("N91").Value = -5121.8
$ws.Range("H98").Value = 7813062
$ws.Range("I98").Value = 7813062
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 7813062
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = -7811564
$ws.Range("N98").ClearContents()
$ws.Range("H116").Value = 14609.772
$ws.Range("I116").Value = 11108.929
$ws.Range("J116").Value = 20736.25
$ws.Range("K116").Value = 11108.929
$ws.Range("L116").Value = 20736.25
$ws.Range("M116").Value = -7666.929
$ws.Range("N116").Value = -27620.25
$ws.Range("H122").Value = 7813062
$ws.Range("I122").Value = 7813062
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 23439186
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -23436736
$ws.Range("N122").ClearContents()
$ws.Range("H124").Value = 86666.336
$ws.Range("J124").Value = 86666.336
$ws.Range("L124").Value = 86666.336
$ws.Range("N124").Value = -96486.336
$ws.Range("H125").Value = 1032
$ws.Range("I125").Value = 1032
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 9288
$ws.Range("L125").Value = 0
$ws.Range("M125").Value = -6828
$ws.Range("N125").ClearContents()
$ws.Range("H127").Value = 2318.9
$ws.Range("I127").Value = 2064.8333
$ws.Range("J127").Value = 2700
$ws.Range("K127").Value = 6194.499899999999
$ws.Range("L127").Value = 8100
$ws.Range("M127").Value = -1234.499899999999
$ws.Range("N127").Value = -18020
$ws.Range("H128").Value = 95999.8
$ws.Range("J128").Value = 95999.8
$ws.Range("L128").Value = 95999.8
$ws.Range("N128").Value = -105959.8
$ws.Range("H129").Value = 2619.4211
$ws.Range("I129").Value = 828.7692
$ws.Range("J129").Value = 6499.1665
$ws.Range("K129").Value = 2486.3076
$ws.Range("L129").Value = 19497.4995
$ws.Range("M129").Value = 2513.6924
$ws.Range("N129").Value = -29497.4995
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()
$ws.Range("H131").Value = 1897130.6
$ws.Range("I131").Value = 2261.125
$ws.Range("K131").Value = 6783.375
$ws.Range("M131").Value = -1743.375
$ws.Range("H141").Value = 3774.3901
$ws.Range("I141").Value = 3067.1316
$ws.Range("K141").Value = 9201.3948
$ws.Range("M141").Value = -4021.3948

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 3289.7
$ws.Range("I88").Value = 2566.6667
$ws.Range("J88").Value = 3599.5715
$ws.Range("K88").Value = 2566.6667
$ws.Range("L88").Value = 3599.5715
$ws.Range("M88").Value = -2160.6667
$ws.Range("N88").Value = -4411.5715
$ws.Range("H91").Value = 3289.7
$ws.Range("I91").Value = 2566.6667
$ws.Range("J91").Value = 3599.5715
$ws.Range("K91").Value = 2566.6667
$ws.Range("L91").Value = 3599.5715
$ws.Range("M91").Value = -1162.6667
$ws.Range("N91").Value = -6407.5715
$ws.Range("H122").Value = 3278.9285
$ws.Range("I122").Value = 3338.923
$ws.Range("K122").Value = 10016.769
$ws.Range("M122").Value = -7566.769
$ws.Range("H123").Value = 329999.34
$ws.Range("J123").Value = 329999.34
$ws.Range("L123").Value = 329999.34
$ws.Range("N123").Value = -339799.34
$ws.Range("H124").Value = 45085.6
$ws.Range("J124").Value = 45085.6
$ws.Range("L124").Value = 45085.6
$ws.Range("N124").Value = -54905.6
$ws.Range("H125").Value = 72300.5
$ws.Range("J125").Value = 72300.5
$ws.Range("L125").Value = 72300.5
$ws.Range("N125").Value = -82140.5
$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()
$ws.Range("H131").Value = 60715
$ws.Range("J131").Value = 60715
$ws.Range("L131").Value = 60715
$ws.Range("N131").Value = -70795

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H123").Value = 10416.6
$ws.Range("I123").Value = 4933.3335
$ws.Range("J123").Value = 18641.5
$ws.Range("K123").Value = 14800.0005
$ws.Range("L123").Value = 55924.5
$ws.Range("M123").Value = -12350.0005
$ws.Range("N123").Value = -60824.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H123").Value = 117999
$ws.Range("J123").Value = 117999
$ws.Range("L123").Value = 117999
$ws.Range("N123").Value = -122899
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()
$ws.Range("H125").Value = 99999.5
$ws.Range("J125").Value = 99999.5
$ws.Range("L125").Value = 99999.5
$ws.Range("N125").Value = -104919.5
$ws.Range("H128").Value = 99999
$ws.Range("J128").Value = 99999
$ws.Range("L128").Value = 99999
$ws.Range("N128").Value = -109959
$ws.Range("H129").Value = 79999.5
$ws.Range("J129").Value = 79999.5
$ws.Range("L129").Value = 79999.5
$ws.Range("N129").Value = -89999.5
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()
$ws.Range("H131").Value = 121497
$ws.Range("J131").Value = 121497
$ws.Range("L131").Value = 121497
$ws.Range("N131").Value = -131577

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H123").Value = 99714
$ws.Range("J123").Value = 99714
$ws.Range("L123").Value = 99714
$ws.Range("N123").Value = -109514
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()
$ws.Range("H131").Value = 99999
$ws.Range("J131").Value = 99999
$ws.Range("L131").Value = 99999
$ws.Range("N131").Value = -110079

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()
$ws.Range("H126").Value = 4095.125
$ws.Range("I126").Value = 2595.7368
$ws.Range("J126").Value = 9792.799999999999
$ws.Range("K126").Value = 7787.2104
$ws.Range("L126").Value = 29378.4
$ws.Range("M126").Value = -5317.2104
$ws.Range("N126").Value = -34318.39999999999
$ws.Range("H127").Value = 42209.5
$ws.Range("I127").Value = 39990
$ws.Range("J127").Value = 44429
$ws.Range("K127").Value = 39990
$ws.Range("L127").Value = 44429
$ws.Range("M127").Value = -35030
$ws.Range("N127").Value = -54349
$ws.Range("H128").Value = 72452.164
$ws.Range("J128").Value = 72452.164
$ws.Range("L128").Value = 72452.164
$ws.Range("N128").Value = -82412.164
$ws.Range("H130").Value = 120000
$ws.Range("J130").Value = 120000
$ws.Range("L130").Value = 120000
$ws.Range("N130").Value = -130040

